$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 60 and 61: the match details in columns F:V were swapped
#    (columns A:E - index/pais/torneio/temporada/data_partida - are
#    unchanged). Use temp holding cells so we don't clobber data we
#    still need to read.
# ------------------------------------------------------------------

# Snapshot current row 60 (F:V) into an unused scratch area (row 1000)
$ws.Range("F60:V60").Copy()
$ws.Range("F1000").PasteSpecial(-4163)

# New row 60 <- old row 61 (F:V)
$ws.Range("F61:V61").Copy()
$ws.Range("F60").PasteSpecial(-4163)

# New row 61 <- old row 60 (F:V), taken from the scratch copy
$ws.Range("F1000:V1000").Copy()
$ws.Range("F61").PasteSpecial(-4163)

# Clean up the scratch row
$ws.Range("F1000:V1000").ClearContents()

# ------------------------------------------------------------------
# 2) Append two new match rows (86 and 87) at the end of the table,
#    copying the formatting of the last existing row (85).
# ------------------------------------------------------------------
$ws.Range("A85:V85").Copy()
$ws.Range("A86:V86").PasteSpecial(-4122)
$ws.Range("A85:V85").Copy()
$ws.Range("A87:V87").PasteSpecial(-4122)

# --- Row 86 ---
$ws.Range("A86").Value2 = 85
$ws.Range("B86").Value2 = "portugal"
$ws.Range("C86").Value2 = "liga-portugal-2"
$ws.Range("D86").Value2 = "2023-2024"
$ws.Range("E86").Value = 45241.6875
$ws.Range("F86").Value2 = "Santa Clara"
$ws.Range("G86").Value2 = 2
$ws.Range("H86").Value2 = "AVS"
$ws.Range("I86").Value2 = 1
$ws.Range("J86").Value2 = 2.13
$ws.Range("K86").Value2 = "08/11/2023 06:12"
$ws.Range("L86").Value2 = 2.19
$ws.Range("M86").Value2 = "11/11/2023 16:22"
$ws.Range("N86").Value2 = 3.31
$ws.Range("O86").Value2 = "08/11/2023 06:12"
$ws.Range("P86").Value2 = 3.17
$ws.Range("Q86").Value2 = "11/11/2023 16:22"
$ws.Range("R86").Value2 = 3.67
$ws.Range("S86").Value2 = "08/11/2023 06:12"
$ws.Range("T86").Value2 = 3.82
$ws.Range("U86").Value2 = "11/11/2023 16:22"
$ws.Range("V86").Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal-2/santa-clara-avs/z1cwL6yo/"

# --- Row 87 ---
$ws.Range("A87").Value2 = 86
$ws.Range("B87").Value2 = "portugal"
$ws.Range("C87").Value2 = "liga-portugal-2"
$ws.Range("D87").Value2 = "2023-2024"
$ws.Range("E87").Value = 45241.79166666666
$ws.Range("F87").Value2 = "Oliveirense"
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = "Pacos Ferreira"
$ws.Range("I87").Value2 = 0
$ws.Range("J87").Value2 = 3.15
$ws.Range("K87").Value2 = "08/11/2023 06:12"
$ws.Range("L87").Value2 = 3.22
$ws.Range("M87").Value2 = "11/11/2023 13:18"
$ws.Range("N87").Value2 = 3.6
$ws.Range("O87").Value2 = "08/11/2023 06:12"
$ws.Range("P87").Value2 = 3.62
$ws.Range("Q87").Value2 = "11/11/2023 13:18"
$ws.Range("R87").Value2 = 2.15
$ws.Range("S87").Value2 = "08/11/2023 06:12"
$ws.Range("T87").Value2 = 2.24
$ws.Range("U87").Value2 = "11/11/2023 13:18"
$ws.Range("V87").Value2 = "https://www.betexplorer.com/football/portugal/liga-portugal-2/oliveirense-pacos-ferreira/d6rYtlb4/"

# ------------------------------------------------------------------
# 3) Refresh the sheet dimension to cover the two new rows.
# ------------------------------------------------------------------
$ws.UsedRange | Out-Null
